$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number (e.g. "213.63").
# Force them to Text format first so Excel keeps the exact original string
# instead of re-interpreting it as a floating point number, then restore the
# default "Normal" style so no stray formatting is left behind on the cell.
$priceTextCells = @("D5", "D8", "D9", "D11", "D16", "D18", "D22", "D25", "D26", "D27", "D28", "D32", "D37", "D39", "D40", "D41", "D44", "D46", "D49")
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.909.60"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "1.647.10"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "213.63"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "23.55"
$ws.Range("E8").Value = "  +3.96%  "
$ws.Range("D9").Value = "0.265"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "0.0872"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").Value = "1.881.64"
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("D13").Value = "1.635.29"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").Value = "65.65"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "27.921.79"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "231.82"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "10.70"
$ws.Range("E22").Value = "  +5.90%  "
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("D25").Value = "151.90"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").Value = "6.92"
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "0.112"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "15.72"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +2.62%  "
$ws.Range("D33").Value = "1.454.66"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").Value = "0.888"
$ws.Range("E37").Value = "  +3.32%  "
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("D39").Value = "0.563"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").Value = "0.920"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").Value = "69.36"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "2.46"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").Value = "5.38"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("E47").Value = "  +5.75%  "
$ws.Range("D48").Value = "1.790.12"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "89.01"
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("E51").Value = "  +1.03%  "

foreach ($addr in $priceTextCells) {
    $ws.Range($addr).Style = "Normal"
}
